# Update row 307 (last existing data row) with revised O/H/L/C values,
# then append 3 new rows (308-310) following the same pattern, copying the
# date-cell formatting from row 307's A cell so new rows match exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 307 ---
$ws.Range("C307").Value = 4812070000000
$ws.Range("D307").Value = 4812070000000
$ws.Range("E307").Value = 4812070000000
$ws.Range("F307").Value = 4812070000000
$ws.Range("G307").Value = 0

# Copy the formatting of the A307 date cell so new date cells (A308:A310)
# pick up the same number format / style as the rest of the column.
$ws.Range("A307").Copy()

# --- New row 308 ---
$ws.Range("A308").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A308").Value = 45108.41666666666
$ws.Range("B308").Value = "ECONOMICS:SEM2"
$ws.Range("C308").Value = 4784709000000
$ws.Range("D308").Value = 4784709000000
$ws.Range("E308").Value = 4784709000000
$ws.Range("F308").Value = 4784709000000
$ws.Range("G308").Value = 0

# --- New row 309 ---
$ws.Range("A309").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A309").Value = 45139.41666666666
$ws.Range("B309").Value = "ECONOMICS:SEM2"
$ws.Range("C309").Value = 4765034000000
$ws.Range("D309").Value = 4765034000000
$ws.Range("E309").Value = 4765034000000
$ws.Range("F309").Value = 4765034000000
$ws.Range("G309").Value = 0

# --- New row 310 ---
$ws.Range("A310").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A310").Value = 45170.41666666666
$ws.Range("B310").Value = "ECONOMICS:SEM2"
$ws.Range("C310").Value = 4744031000000
$ws.Range("D310").Value = 4744031000000
$ws.Range("E310").Value = 4744031000000
$ws.Range("F310").Value = 4744031000000
$ws.Range("G310").Value = 0

$excel.CutCopyMode = $false
